$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.819.55"
$ws.Range("E2").Value = "  -7.01%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.219.10"
$ws.Range("E3").Value = "  -8.36%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "510.19"
$ws.Range("E5").Value = "  -7.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.73"
$ws.Range("E6").Value = "  -15.32%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.586"
$ws.Range("E7").Value = "  -7.08%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.217.58"
$ws.Range("E9").Value = "  -8.20%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.588"
$ws.Range("E10").Value = "  -10.02%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.95"
$ws.Range("E11").Value = "  -11.75%  "
$ws.Range("E12").Value = "  -9.36%  "
$ws.Range("E13").Value = "  -6.89%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.79"
$ws.Range("E14").Value = "  -10.22%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.742.25"
$ws.Range("E15").Value = "  -8.12%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.225.23"
$ws.Range("E16").Value = "  -8.26%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.704.50"
$ws.Range("E17").Value = "  -6.79%  "
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.113"
$ws.Range("E18").Value = "  -9.17%  "
$ws.Range("E19").Value = "  -7.88%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.67"
$ws.Range("E20").Value = "  -9.59%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.932"
$ws.Range("E21").Value = "  -8.98%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "362.00"
$ws.Range("E22").Value = "  -8.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.64"
$ws.Range("E23").Value = "  -8.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.28"
$ws.Range("E24").Value = "  -7.96%  "
$ws.Range("B25").Value = "LEO"
$ws.Range("C25").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "6.15"
$ws.Range("E25").Value = "  -0.27%  "
$ws.Range("B26").Value = "RenderToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.60"
$ws.Range("E26").Value = "  -11.16%  "
$ws.Range("E27").Value = "  -1.09%  "
$ws.Range("E28").Value = "  -7.77%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.98"
$ws.Range("E29").Value = "  -10.61%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.05"
$ws.Range("E30").Value = "  -9.10%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "27.94"
$ws.Range("E31").Value = "  -10.66%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "620.88"
$ws.Range("E32").Value = "  -13.46%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.38"
$ws.Range("E33").Value = "  -8.81%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "10.93"
$ws.Range("E34").Value = "  -6.49%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "57.53"
$ws.Range("E35").Value = "  -10.07%  "
$ws.Range("E36").Value = "  -7.93%  "
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "35.79"
$ws.Range("E38").Value = "  -6.66%  "
$ws.Range("E39").Value = "  -5.19%  "
$ws.Range("E40").Value = "  +0.11%  "
$ws.Range("E41").Value = "  -2.13%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.842.07"
$ws.Range("E42").Value = "  -7.08%  "
$ws.Range("E43").Value = "  -9.16%  "
$ws.Range("B44").Value = "Fetch.AI"
$ws.Range("C44").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.37"
$ws.Range("E44").Value = "  -5.52%  "
$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.58"
$ws.Range("E45").Value = "  -7.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.58"
$ws.Range("E46").Value = "  -14.32%  "
$ws.Range("E47").Value = "  -6.43%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.96"
$ws.Range("E48").Value = "  +0.85%  "
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.69"
$ws.Range("E49").Value = "  +2.53%  "
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "133.26"
$ws.Range("E50").Value = "  -3.88%  "
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.121"
$ws.Range("E51").Value = "  -6.93%  "
